$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'29.279.46"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$c = $ws.Range("D3")
$c.Value = "'1.863.98"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "
$c = $ws.Range("D4")
$c.Value = "'0.9995"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.Value = "'0.7092"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$c = $ws.Range("D6")
$c.Value = "'237.76"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$c = $ws.Range("D7")
$c.Value = "'0.9998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.Value = "'0.07899"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.94%  "
$c = $ws.Range("D9")
$c.Value = "'0.3061"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.43%  "
$c = $ws.Range("D10")
$c.Value = "'25.00"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.99%  "
$c = $ws.Range("D11")
$c.Value = "'0.08152"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.37%  "
$c = $ws.Range("D12")
$c.Value = "'1.862.52"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "
$c = $ws.Range("D13")
$c.Value = "'5.220"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "
$c = $ws.Range("D14")
$c.Value = "'0.7193"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "
$c = $ws.Range("D15")
$c.Value = "'89.09"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.31%  "
$c = $ws.Range("D16")
$c.Value = "'29.437.38"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$c = $ws.Range("D17")
$c.Value = "'5.803"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$c = $ws.Range("D18")
$c.Value = "'241.62"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "
$c = $ws.Range("D19")
$c.Value = "'0.000007797"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$c = $ws.Range("D20")
$c.Value = "'13.20"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D21")
$c.Value = "'2.143.69"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D22")
$c.Value = "'0.9993"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$c = $ws.Range("D23")
$c.Value = "'0.9991"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$c = $ws.Range("D24")
$c.Value = "'7.581"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "
$c = $ws.Range("D25")
$c.Value = "'162.31"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D26")
$c.Value = "'8.913"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D27")
$c.Value = "'0.1447"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$c = $ws.Range("D28")
$c.Value = "'18.09"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$c = $ws.Range("D29")
$c.Value = "'1.909"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "
$c = $ws.Range("D30")
$c.Value = "'1.368"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.68%  "
$c = $ws.Range("D31")
$c.Value = "'1.473"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$c = $ws.Range("D32")
$c.Value = "'4.312"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.58%  "
$c = $ws.Range("D33")
$c.Value = "'4.047"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "
$c = $ws.Range("D34")
$c.Value = "'0.05184"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "
$c = $ws.Range("D35")
$c.Value = "'1.183"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "
$c = $ws.Range("D36")
$c.Value = "'0.7176"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.44%  "
$c = $ws.Range("D37")
$c.Value = "'1.007"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
$c = $ws.Range("D38")
$c.Value = "'2.674"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$c = $ws.Range("D39")
$c.Value = "'0.01849"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$c = $ws.Range("D40")
$c.Value = "'2.695"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "
$c = $ws.Range("D41")
$c.Value = "'1.167.86"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "
$c = $ws.Range("D42")
$c.Value = "'0.9159"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "
$c = $ws.Range("D43")
$c.Value = "'5.998"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$c = $ws.Range("D44")
$c.Value = "'71.59"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$c = $ws.Range("D45")
$c.Value = "'0.4268"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "
$c = $ws.Range("D46")
$c.Value = "'0.9993"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$c = $ws.Range("D47")
$c.Value = "'102.08"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "
$c = $ws.Range("D48")
$c.Value = "'0.5346"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "
$c = $ws.Range("D49")
$c.Value = "'1.753"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "
$c = $ws.Range("D50")
$c.Value = "'9.167"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$c = $ws.Range("D51")
$c.Value = "'6.973"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
